$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "42.941.86"
Set-TextValue $ws "E2" "  -1.55%  "
Set-TextValue $ws "D3" "2.297.46"
Set-TextValue $ws "E3" "  -1.92%  "
Set-TextValue $ws "E4" "  +0.05%  "
Set-TextValue $ws "D5" "300.17"
Set-TextValue $ws "E5" "  -1.63%  "
Set-TextValue $ws "D6" "97.02"
Set-TextValue $ws "E6" "  -5.04%  "
Set-TextValue $ws "D7" "0.505"
Set-TextValue $ws "E7" "  -1.47%  "
Set-TextValue $ws "E8" "  +0.03%  "
Set-TextValue $ws "D9" "0.495"
Set-TextValue $ws "E9" "  -4.72%  "
Set-TextValue $ws "D10" "33.57"
Set-TextValue $ws "E10" "  -5.05%  "
Set-TextValue $ws "D11" "0.0794"
Set-TextValue $ws "E11" "  -0.43%  "
Set-TextValue $ws "D12" "49.30"
Set-TextValue $ws "E12" "  -4.77%  "
Set-TextValue $ws "E13" "  +1.97%  "
Set-TextValue $ws "D14" "16.77"
Set-TextValue $ws "E14" "  +7.53%  "
Set-TextValue $ws "D15" "6.76"
Set-TextValue $ws "E15" "  -1.10%  "
Set-TextValue $ws "D16" "2.657.16"
Set-TextValue $ws "E16" "  -1.83%  "
Set-TextValue $ws "D17" "2.316.85"
Set-TextValue $ws "E17" "  -0.99%  "
Set-TextValue $ws "D18" "0.809"
Set-TextValue $ws "E18" "  -0.20%  "
Set-TextValue $ws "D19" "42.847.58"
Set-TextValue $ws "E19" "  -1.55%  "
Set-TextValue $ws "B20" "ShibaInu"
Set-TextValue $ws "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D20" "0.0₃0899"
Set-TextValue $ws "E20" "  -1.07%  "
Set-TextValue $ws "B21" "InternetComputer(DFINITY)"
Set-TextValue $ws "C21" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D21" "11.57"
Set-TextValue $ws "E21" "  -2.13%  "
Set-TextValue $ws "E22" "  -2.25%  "
Set-TextValue $ws "D23" "67.22"
Set-TextValue $ws "E23" "  -1.35%  "
Set-TextValue $ws "D24" "235.90"
Set-TextValue $ws "E24" "  -1.21%  "
Set-TextValue $ws "E25" "  +0.87%  "
Set-TextValue $ws "E26" "  +0.03%  "
Set-TextValue $ws "D27" "2.45"
Set-TextValue $ws "E27" "  -3.78%  "
Set-TextValue $ws "D28" "24.32"
Set-TextValue $ws "E28" "  -3.13%  "
Set-TextValue $ws "E29" "  -6.16%  "
Set-TextValue $ws "D30" "166.51"
Set-TextValue $ws "E30" "  +0.61%  "
Set-TextValue $ws "D31" "33.82"
Set-TextValue $ws "E31" "  -2.72%  "
Set-TextValue $ws "D32" "9.10"
Set-TextValue $ws "E32" "  -2.06%  "
Set-TextValue $ws "E33" "  +0.07%  "
Set-TextValue $ws "D34" "4.72"
Set-TextValue $ws "E34" "  +3.99%  "
Set-TextValue $ws "D35" "4.94"
Set-TextValue $ws "E35" "  -2.82%  "
Set-TextValue $ws "E36" "  -1.51%  "
Set-TextValue $ws "D37" "16.80"
Set-TextValue $ws "E37" "  -1.82%  "
Set-TextValue $ws "D38" "0.0691"
Set-TextValue $ws "E38" "  -2.54%  "
Set-TextValue $ws "B39" "LidoDAOToken"
Set-TextValue $ws "C39" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D39" "2.81"
Set-TextValue $ws "E39" "  -3.23%  "
Set-TextValue $ws "B40" "Kaspa"
Set-TextValue $ws "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D40" "0.101"
Set-TextValue $ws "E40" "  -2.04%  "
Set-TextValue $ws "D41" "1.74"
Set-TextValue $ws "E41" "  -5.00%  "
Set-TextValue $ws "E42" "  -2.33%  "
Set-TextValue $ws "E43" "  -3.10%  "
Set-TextValue $ws "D44" "1.982.97"
Set-TextValue $ws "E44" "  -0.42%  "
Set-TextValue $ws "D45" "0.0280"
Set-TextValue $ws "E45" "  -1.95%  "
Set-TextValue $ws "D46" "17.67"
Set-TextValue $ws "D47" "9.81"
Set-TextValue $ws "E47" "  -1.51%  "
Set-TextValue $ws "D48" "2.84"
Set-TextValue $ws "E48" "  -3.99%  "
Set-TextValue $ws "D49" "2.525.84"
Set-TextValue $ws "E49" "  -1.64%  "
Set-TextValue $ws "D50" "52.74"
Set-TextValue $ws "E50" "  -6.49%  "
Set-TextValue $ws "D51" "4.56"
Set-TextValue $ws "E51" "  -7.18%  "
